# Refresh the crypto price / 1h-volume snapshot and rotate the
# GateToken..LEO symbol block (rows 6-18) per the upstream feed pull
# that produced "Updated symbol list on Tue Jan 31 20:31:00 UTC 2023".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.70"
$ws.Range("E2").Value = "'2.31%"
$ws.Range("D3").Value = "'37.75"
$ws.Range("E3").Value = "'1.91%"
$ws.Range("D4").Value = "'5.144"
$ws.Range("E4").Value = "'0.78%"
$ws.Range("D5").Value = "'0.07899"
$ws.Range("E5").Value = "'2.44%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.915"
$ws.Range("E6").Value = "'2.48%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.317"
$ws.Range("E7").Value = "'1.32%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.862"
$ws.Range("E8").Value = "'-9.78%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9221"
$ws.Range("E9").Value = "'0.30%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1238"
$ws.Range("E10").Value = "'0.86%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1928"
$ws.Range("E11").Value = "'1.90%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09166"
$ws.Range("E12").Value = "'5.34%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03345"
$ws.Range("E13").Value = "'-1.75%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09604"
$ws.Range("E14").Value = "'-0.90%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001384"
$ws.Range("E15").Value = "'1.31%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005830"
$ws.Range("E16").Value = "'-4.33%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.503"
$ws.Range("E17").Value = "'-1.48%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.425"
$ws.Range("E18").Value = "'1.05%"
$ws.Range("D19").Value = "'0.3442"
$ws.Range("E19").Value = "'2.05%"
$ws.Range("D20").Value = "'5.275"
$ws.Range("E20").Value = "'4.88%"
$ws.Range("E21").Value = "'-0.86%"
$ws.Range("D22").Value = "'0.2591"
$ws.Range("E22").Value = "'3.51%"
$ws.Range("E23").Value = "'-0.63%"
$ws.Range("D24").Value = "'0.04381"
$ws.Range("E24").Value = "'1.20%"
$ws.Range("D25").Value = "'0.001249"
$ws.Range("E25").Value = "'2.42%"
$ws.Range("E26").Value = "'-3.26%"
$ws.Range("D27").Value = "'0.0001220"
$ws.Range("E27").Value = "'-10.01%"
$ws.Range("D39").Value = "'0.02281"
$ws.Range("E39").Value = "'2.51%"
$ws.Range("D40").Value = "'0.05099"
$ws.Range("E40").Value = "'4.04%"
$ws.Range("D41").Value = "'0.007458"
$ws.Range("E41").Value = "'-2.17%"
$ws.Range("D42").Value = "'0.1360"
$ws.Range("E42").Value = "'2.08%"
$ws.Range("D43").Value = "'0.008800"
$ws.Range("E43").Value = "'-11.38%"
$ws.Range("D44").Value = "'0.001961"
$ws.Range("E44").Value = "'-2.11%"
$ws.Range("D45").Value = "'0.008613"
$ws.Range("E45").Value = "'-2.31%"
$ws.Range("D46").Value = "'0.00006743"
$ws.Range("E46").Value = "'-3.20%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.56%"
$ws.Range("E48").Value = "'11.20%"
$ws.Range("D49").Value = "'0.001200"
$ws.Range("E49").Value = "'-8.21%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.56%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.56%"
